$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (current outputTemperaturePython/FMU columns)
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# Copy the header formatting from the (now shifted) outputTemperaturePython header into the new cells
$ws.Range("F1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Set the new header labels
$ws.Range("D1").Value = "controlSignalPython"
$ws.Range("E1").Value = "controlSignalFMU"

# Fill in the new column values
$ws.Range("D2").Value = 2921
$ws.Range("E2").Value = 2921

$ws.Range("D3").Value = 0.696351587164061
$ws.Range("E3").Value = 0.6961559566308066

$ws.Range("D4").Value = 0.4168086750774854
$ws.Range("E4").Value = 0.4170013826581315

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.6497365006255256
$ws.Range("E6").Value = 0.6497365006255256

$ws.Range("D7").Value = 2.174568593007422
$ws.Range("E7").Value = 2.174568593007422
